$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginPage")

# Fill in the two newly added open questions in column C (rows 7 and 8)
$ws.Range("C7").Value = "How is OTP Generated and validated."
$ws.Range("C8").Value = "How is new Password Generated"

# Move the active selection to C8 as reflected in the saved view state
$ws.Activate()
$ws.Range("C8").Select()
